$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 7).Value2 = 2.441378666666667
$ws.Cells.Item(2, 8).Value2 = 7.324135999999999
$ws.Cells.Item(2, 9).Value2 = 0.1119936059016048
$ws.Cells.Item(2, 10).Value2 = 0.1119936059016048
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 13).Value2 = 233.0249173333333
$ws.Cells.Item(2, 14).Value2 = 699.074752
$ws.Cells.Item(2, 15).Value2 = 0.765932814443446
$ws.Cells.Item(2, 16).Value2 = 0.765932814443446
$ws.Cells.Item(2, 17).Value2 = 568.9020619793635
$ws.Cells.Item(2, 18).Value2 = 5120.118557814271
$ws.Cells.Item(2, 19).Value2 = 0.08577957776788628
$ws.Cells.Item(2, 20).Value2 = 0.08577957776788628

$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 7).Value2 = 2.441378666666667
$ws.Cells.Item(3, 8).Value2 = 7.324135999999999
$ws.Cells.Item(3, 9).Value2 = 0.1119936059016048
$ws.Cells.Item(3, 10).Value2 = 0.1119936059016048
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 13).Value2 = 51.70670833333333
$ws.Cells.Item(3, 14).Value2 = 155.120125
$ws.Cells.Item(3, 15).Value2 = 0.1699554927111273
$ws.Cells.Item(3, 16).Value2 = 0.1699554927111273
$ws.Cells.Item(3, 17).Value2 = 126.2356546485556
$ws.Cells.Item(3, 18).Value2 = 1136.120891837
$ws.Cells.Item(3, 19).Value2 = 0.01903392847150306
$ws.Cells.Item(3, 20).Value2 = 0.01903392847150306

$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 7).Value2 = 2.441378666666667
$ws.Cells.Item(4, 8).Value2 = 7.324135999999999
$ws.Cells.Item(4, 9).Value2 = 0.1119936059016048
$ws.Cells.Item(4, 10).Value2 = 0.1119936059016048
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 13).Value2 = 8.768542333333334
$ws.Cells.Item(4, 14).Value2 = 26.305627
$ws.Cells.Item(4, 15).Value2 = 0.02882144272292286
$ws.Cells.Item(4, 16).Value2 = 0.02882144272292286
$ws.Cells.Item(4, 17).Value2 = 21.40733219036356
$ws.Cells.Item(4, 18).Value2 = 192.665989713272
$ws.Cells.Item(4, 19).Value2 = 0.003227817297826699
$ws.Cells.Item(4, 20).Value2 = 0.003227817297826698

$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 7).Value2 = 2.441378666666667
$ws.Cells.Item(5, 8).Value2 = 7.324135999999999
$ws.Cells.Item(5, 9).Value2 = 0.1119936059016048
$ws.Cells.Item(5, 10).Value2 = 0.1119936059016048
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 13).Value2 = 6.761708
$ws.Cells.Item(5, 14).Value2 = 20.285124
$ws.Cells.Item(5, 15).Value2 = 0.02222515127631772
$ws.Cells.Item(5, 16).Value2 = 0.02222515127631772
$ws.Cells.Item(5, 17).Value2 = 16.50788966142933
$ws.Cells.Item(5, 18).Value2 = 148.571006952864
$ws.Cells.Item(5, 19).Value2 = 0.002489074833143475
$ws.Cells.Item(5, 20).Value2 = 0.002489074833143475

$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 7).Value2 = 2.441378666666667
$ws.Cells.Item(6, 8).Value2 = 7.324135999999999
$ws.Cells.Item(6, 9).Value2 = 0.1119936059016048
$ws.Cells.Item(6, 10).Value2 = 0.1119936059016048
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 13).Value2 = 3.974883333333333
$ws.Cells.Item(6, 14).Value2 = 11.92465
$ws.Cells.Item(6, 15).Value2 = 0.0130650988461861
$ws.Cells.Item(6, 16).Value2 = 0.01306509884618611
$ws.Cells.Item(6, 17).Value2 = 9.704195372488888
$ws.Cells.Item(6, 18).Value2 = 87.33775835239999
$ws.Cells.Item(6, 19).Value2 = 0.001463207531245278
$ws.Cells.Item(6, 20).Value2 = 0.001463207531245278

$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 7).Value2 = 12.88577866666667
$ws.Cells.Item(7, 8).Value2 = 38.657336
$ws.Cells.Item(7, 9).Value2 = 0.5911106037886134
$ws.Cells.Item(7, 10).Value2 = 0.5911106037886134
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 13).Value2 = 233.0249173333333
$ws.Cells.Item(7, 14).Value2 = 699.074752
$ws.Cells.Item(7, 15).Value2 = 0.765932814443446
$ws.Cells.Item(7, 16).Value2 = 0.765932814443446
$ws.Cells.Item(7, 17).Value2 = 3002.70750857563
$ws.Cells.Item(7, 18).Value2 = 27024.36757718067
$ws.Cells.Item(7, 19).Value2 = 0.4527510084071773
$ws.Cells.Item(7, 20).Value2 = 0.4527510084071773

$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 7).Value2 = 12.88577866666667
$ws.Cells.Item(8, 8).Value2 = 38.657336
$ws.Cells.Item(8, 9).Value2 = 0.5911106037886134
$ws.Cells.Item(8, 10).Value2 = 0.5911106037886134
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 13).Value2 = 51.70670833333333
$ws.Cells.Item(8, 14).Value2 = 155.120125
$ws.Cells.Item(8, 15).Value2 = 0.1699554927111273
$ws.Cells.Item(8, 16).Value2 = 0.1699554927111273
$ws.Cells.Item(8, 17).Value2 = 666.2811991652222
$ws.Cells.Item(8, 18).Value2 = 5996.530792487
$ws.Cells.Item(8, 19).Value2 = 0.1004624939136657
$ws.Cells.Item(8, 20).Value2 = 0.1004624939136658

$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 7).Value2 = 12.88577866666667
$ws.Cells.Item(9, 8).Value2 = 38.657336
$ws.Cells.Item(9, 9).Value2 = 0.5911106037886134
$ws.Cells.Item(9, 10).Value2 = 0.5911106037886134
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 13).Value2 = 8.768542333333334
$ws.Cells.Item(9, 14).Value2 = 26.305627
$ws.Cells.Item(9, 15).Value2 = 0.02882144272292286
$ws.Cells.Item(9, 16).Value2 = 0.02882144272292286
$ws.Cells.Item(9, 17).Value2 = 112.9894957366302
$ws.Cells.Item(9, 18).Value2 = 1016.905461629672
$ws.Cells.Item(9, 19).Value2 = 0.01703666041000587
$ws.Cells.Item(9, 20).Value2 = 0.01703666041000587

$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 7).Value2 = 12.88577866666667
$ws.Cells.Item(10, 8).Value2 = 38.657336
$ws.Cells.Item(10, 9).Value2 = 0.5911106037886134
$ws.Cells.Item(10, 10).Value2 = 0.5911106037886134
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 13).Value2 = 6.761708
$ws.Cells.Item(10, 14).Value2 = 20.285124
$ws.Cells.Item(10, 15).Value2 = 0.02222515127631772
$ws.Cells.Item(10, 16).Value2 = 0.02222515127631772
$ws.Cells.Item(10, 17).Value2 = 87.12987269662933
$ws.Cells.Item(10, 18).Value2 = 784.168854269664
$ws.Cells.Item(10, 19).Value2 = 0.01313752259023744
$ws.Cells.Item(10, 20).Value2 = 0.01313752259023744

$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 7).Value2 = 12.88577866666667
$ws.Cells.Item(11, 8).Value2 = 38.657336
$ws.Cells.Item(11, 9).Value2 = 0.5911106037886134
$ws.Cells.Item(11, 10).Value2 = 0.5911106037886134
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 13).Value2 = 3.974883333333333
$ws.Cells.Item(11, 14).Value2 = 11.92465
$ws.Cells.Item(11, 15).Value2 = 0.0130650988461861
$ws.Cells.Item(11, 16).Value2 = 0.01306509884618611
$ws.Cells.Item(11, 17).Value2 = 51.21946685915555
$ws.Cells.Item(11, 18).Value2 = 460.9752017324
$ws.Cells.Item(11, 19).Value2 = 0.007722918467526984
$ws.Cells.Item(11, 20).Value2 = 0.007722918467526985

$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 7).Value2 = 6.472111000000001
$ws.Cells.Item(12, 8).Value2 = 19.416333
$ws.Cells.Item(12, 9).Value2 = 0.2968957903097819
$ws.Cells.Item(12, 10).Value2 = 0.2968957903097818
$ws.Cells.Item(12, 11).Value2 = 3
$ws.Cells.Item(12, 13).Value2 = 233.0249173333333
$ws.Cells.Item(12, 14).Value2 = 699.074752
$ws.Cells.Item(12, 15).Value2 = 0.765932814443446
$ws.Cells.Item(12, 16).Value2 = 0.765932814443446
$ws.Cells.Item(12, 17).Value2 = 1508.163130747158
$ws.Cells.Item(12, 18).Value2 = 13573.46817672442
$ws.Cells.Item(12, 19).Value2 = 0.2274022282683824
$ws.Cells.Item(12, 20).Value2 = 0.2274022282683824

$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 7).Value2 = 6.472111000000001
$ws.Cells.Item(13, 8).Value2 = 19.416333
$ws.Cells.Item(13, 9).Value2 = 0.2968957903097819
$ws.Cells.Item(13, 10).Value2 = 0.2968957903097818
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 13).Value2 = 51.70670833333333
$ws.Cells.Item(13, 14).Value2 = 155.120125
$ws.Cells.Item(13, 15).Value2 = 0.1699554927111273
$ws.Cells.Item(13, 16).Value2 = 0.1699554927111273
$ws.Cells.Item(13, 17).Value2 = 334.6515557779584
$ws.Cells.Item(13, 18).Value2 = 3011.864002001625
$ws.Cells.Item(13, 19).Value2 = 0.05045907032595851
$ws.Cells.Item(13, 20).Value2 = 0.05045907032595851

$ws.Cells.Item(14, 5).Value2 = 3
$ws.Cells.Item(14, 7).Value2 = 6.472111000000001
$ws.Cells.Item(14, 8).Value2 = 19.416333
$ws.Cells.Item(14, 9).Value2 = 0.2968957903097819
$ws.Cells.Item(14, 10).Value2 = 0.2968957903097818
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 13).Value2 = 8.768542333333334
$ws.Cells.Item(14, 14).Value2 = 26.305627
$ws.Cells.Item(14, 15).Value2 = 0.02882144272292286
$ws.Cells.Item(14, 16).Value2 = 0.02882144272292286
$ws.Cells.Item(14, 17).Value2 = 56.75097928953235
$ws.Cells.Item(14, 18).Value2 = 510.7588136057911
$ws.Cells.Item(14, 19).Value2 = 0.008556965015090295
$ws.Cells.Item(14, 20).Value2 = 0.008556965015090292

$ws.Cells.Item(15, 5).Value2 = 3
$ws.Cells.Item(15, 7).Value2 = 6.472111000000001
$ws.Cells.Item(15, 8).Value2 = 19.416333
$ws.Cells.Item(15, 9).Value2 = 0.2968957903097819
$ws.Cells.Item(15, 10).Value2 = 0.2968957903097818
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 13).Value2 = 6.761708
$ws.Cells.Item(15, 14).Value2 = 20.285124
$ws.Cells.Item(15, 15).Value2 = 0.02222515127631772
$ws.Cells.Item(15, 16).Value2 = 0.02222515127631772
$ws.Cells.Item(15, 17).Value2 = 43.762524725588
$ws.Cells.Item(15, 18).Value2 = 393.862722530292
$ws.Cells.Item(15, 19).Value2 = 0.006598553852936805
$ws.Cells.Item(15, 20).Value2 = 0.006598553852936805

$ws.Cells.Item(16, 5).Value2 = 3
$ws.Cells.Item(16, 7).Value2 = 6.472111000000001
$ws.Cells.Item(16, 8).Value2 = 19.416333
$ws.Cells.Item(16, 9).Value2 = 0.2968957903097819
$ws.Cells.Item(16, 10).Value2 = 0.2968957903097818
$ws.Cells.Item(16, 11).Value2 = 3
$ws.Cells.Item(16, 13).Value2 = 3.974883333333333
$ws.Cells.Item(16, 14).Value2 = 11.92465
$ws.Cells.Item(16, 15).Value2 = 0.0130650988461861
$ws.Cells.Item(16, 16).Value2 = 0.01306509884618611
$ws.Cells.Item(16, 17).Value2 = 25.72588614538333
$ws.Cells.Item(16, 18).Value2 = 231.53297530845
$ws.Cells.Item(16, 19).Value2 = 0.003878972847413843
$ws.Cells.Item(16, 20).Value2 = 0.003878972847413842
